$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 3021
$ws.Cells.Item(7, 6).Value = 2343
$ws.Cells.Item(8, 6).Value = 1736
$ws.Cells.Item(9, 6).Value = 1736
$ws.Cells.Item(10, 6).Value = 64
$ws.Cells.Item(12, 6).Value = 142
$ws.Cells.Item(15, 6).Value = 2685
$ws.Cells.Item(17, 6).Value = 1557
$ws.Cells.Item(18, 6).Value = 7194
$ws.Cells.Item(20, 6).Value = 7340
$ws.Cells.Item(21, 6).Value = 6
$ws.Cells.Item(22, 6).Value = 14
$ws.Cells.Item(23, 6).Value = 5679
$ws.Cells.Item(24, 6).Value = 5679
$ws.Cells.Item(25, 6).Value = 3148
$ws.Cells.Item(26, 6).Value = 3519
$ws.Cells.Item(28, 6).Value = 9
$ws.Cells.Item(31, 6).Value = 1957
$ws.Cells.Item(33, 6).Value = 315
$ws.Cells.Item(35, 6).Value = 233
$ws.Cells.Item(36, 6).Value = 502
$ws.Cells.Item(38, 6).Value = 2481
$ws.Cells.Item(39, 6).Value = 1284
$ws.Cells.Item(40, 6).Value = 2893
$ws.Cells.Item(41, 6).Value = 74
$ws.Cells.Item(43, 6).Value = 178
$ws.Cells.Item(44, 6).Value = 415
$ws.Cells.Item(45, 6).Value = 1128
$ws.Cells.Item(47, 6).Value = 492
$ws.Cells.Item(48, 6).Value = 545

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 9
$ws.Cells.Item(7, 6).Value = 42
$ws.Cells.Item(8, 6).Value = 226
$ws.Cells.Item(12, 6).Value = 362
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(21, 6).Value = 9

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 9
$ws.Cells.Item(5, 6).Value = 3021
$ws.Cells.Item(6, 6).Value = 2343
$ws.Cells.Item(7, 6).Value = 1736
$ws.Cells.Item(8, 6).Value = 1736
$ws.Cells.Item(10, 6).Value = 64
$ws.Cells.Item(12, 6).Value = 142
$ws.Cells.Item(14, 6).Value = 42
$ws.Cells.Item(16, 6).Value = 2685
$ws.Cells.Item(17, 6).Value = 1557
$ws.Cells.Item(20, 6).Value = 7194
$ws.Cells.Item(22, 6).Value = 7340
$ws.Cells.Item(23, 6).Value = 14
$ws.Cells.Item(24, 6).Value = 5679
$ws.Cells.Item(25, 6).Value = 5679
$ws.Cells.Item(26, 6).Value = 3148
$ws.Cells.Item(27, 6).Value = 3519
$ws.Cells.Item(29, 6).Value = 9
$ws.Cells.Item(32, 6).Value = 1957
$ws.Cells.Item(35, 6).Value = 315
$ws.Cells.Item(37, 6).Value = 502
$ws.Cells.Item(39, 6).Value = 2481
$ws.Cells.Item(40, 6).Value = 1284
$ws.Cells.Item(42, 6).Value = 2893
$ws.Cells.Item(43, 6).Value = 74
$ws.Cells.Item(45, 6).Value = 178
$ws.Cells.Item(46, 6).Value = 9
$ws.Cells.Item(47, 6).Value = 415
$ws.Cells.Item(48, 6).Value = 1128
$ws.Cells.Item(50, 6).Value = 492
$ws.Cells.Item(51, 6).Value = 545
